$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Find every cell in column G that contains the exact "Recorded By" value
# "System, dnasr281@gmail.com" and flip the order of the two names to
# "dnasr281@gmail.com, System".
$colG = $ws.Range("G1:G235")

foreach ($cell in $colG.Cells) {
    if ($cell.Text -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
